# Apply the Trace_Report_WCS data update to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: updated search-completed timestamp in the description text
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:10:13 EDT, by WPJTOWN1.The search returned: 8 events."

# Row 9 (car MWCX102553): location/time/event updated
$ws.Range("C9").Value = "DENVER"
$ws.Range("G9").Value = 1323
$ws.Range("H9").Value = "Arrive In-Transit"

# Row 10 (car MWCX100715): location/state/time/event updated
$ws.Range("C10").Value = "BIRMINGHAM"
$ws.Range("D10").Value = "AL"
$ws.Range("G10").Value = 1610
$ws.Range("H10").Value = "Arrive In-Transit"

# Row 11 (car MWCX100705): location/day/time/train id updated
$ws.Range("C11").Value = "NETTLETON"
$ws.Range("F11").Value = 21
$ws.Range("G11").Value = 2000
$ws.Range("I11").Value = "YAMO10"

# Column C (Location City) re-sized to fit the new, shorter city names
$ws.Columns.Item(3).ColumnWidth = 13.7
